# Applies the "role.docx" documentation fix:
#  1. Trim the redundant "AMČR, " prefix (and stray comma) from the
#     "který zavítá do ..." sentence in the first paragraph.
#  2. Extend the "Rozpracované záznamy ..." sentence with the new
#     "uživatelům ze stejné organizace" clause.
#  3. Un-bold the role-name cells ("Anonym", "Badatel", "Archeolog",
#     "Archivář", "Administrátor") in the second column of the role table.
#  4. Nudge the table's column widths to their new values.

$d = $word.ActiveDocument

# --- 1. "AMČR, Digitálního archivu AMČR, či" -> "Digitálního archivu AMČR či"
# The target text immediately follows a bold run ("Anonym"); starting the
# search exactly on the run boundary makes the engine fold the replacement
# into that bold run, so the search is anchored one character past the
# boundary (into the un-bold run) to keep "Anonym" separately bold.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute(") z", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterAnonym = $d.Range($anchor.Start + 1, $d.Content.End)
$afterAnonym.Find.Execute(
    "získává každý uživatel, který zavítá do AMČR, Digitálního archivu AMČR, či jiné aplikace AIS CR.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "získává každý uživatel, který zavítá do Digitálního archivu AMČR či jiné aplikace AIS CR.",
    2) | Out-Null

# --- 2. "který jej edituje, a uživatelům ..." -> "který je edituje, uživatelům ze stejné organizace a uživatelům ..."
$d.Content.Find.Execute(
    "Rozpracované záznamy jsou zpřístupňovány jen tomu uživateli, který jej edituje, a uživatelům s vyššími oprávněními (úroveň D a E).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Rozpracované záznamy jsou zpřístupňovány jen tomu uživateli, který je edituje, uživatelům ze stejné organizace a uživatelům s vyššími oprávněními (úroveň D a E).",
    2) | Out-Null

# --- 3. Remove bold from the role-name cells in the table's "Název" column.
$t = $d.Tables.Item(1)
$roleNames = @("Anonym", "Badatel", "Archeolog", "Archivář", "Administrátor")
for ($row = 2; $row -le $t.Rows.Count; $row++) {
    $cell = $t.Cell($row, 2)
    if ($roleNames -contains $cell.Range.Text.Trim()) {
        $cell.Range.Font.Bold = 0
    }
}

# --- 4. Update the column widths (values given in points; Word stores dxa == pt*20).
$t.Columns.Item(1).Width = 12.5
$t.Columns.Item(2).Width = 29.15
$t.Columns.Item(3).Width = 164.65
$t.Columns.Item(4).Width = 189.65

Write-Host "Edits applied."
